# Weekly fruit/vegetable data update: insert a new row of data
# (row 322) for "Plátano" at "Vega Monumental Concepción", pushing the
# existing rows 322:400 down to 323:401.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 322, shifting everything
# below (rows 322-400) down by one (to 323-401).
$ws.Rows.Item(322).Insert()

# Populate the newly inserted row with this week's entry.
$ws.Cells.Item(322, 1).Value = 11
$ws.Cells.Item(322, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(322, 3).Value = "Bíobío"
$ws.Cells.Item(322, 4).Value = "12/21/2021"
$ws.Cells.Item(322, 5).Value = 8
$ws.Cells.Item(322, 6).Value = "Fruta"
$ws.Cells.Item(322, 7).Value = 100108
$ws.Cells.Item(322, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(322, 9).Value = 100108006
$ws.Cells.Item(322, 10).Value = "Plátano"
$ws.Cells.Item(322, 11).Value = "Sin especificar"
$ws.Cells.Item(322, 12).Value = "Pintón"
$ws.Cells.Item(322, 13).Value = 650
$ws.Cells.Item(322, 14).Value = 10000
$ws.Cells.Item(322, 15).Value = 12000
$ws.Cells.Item(322, 16).Value = 11015
$ws.Cells.Item(322, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(322, 18).Value = "Ecuador"
$ws.Cells.Item(322, 19).Value = 551
$ws.Cells.Item(322, 20).Value = 20
